# feat: add 2022-Q4 data
#
# The workbook only has "总计" (summary) and "2022-Q3" (fund detail) sheets.
# This change adds a new "2022-Q4" fund-detail sheet (placed between the two
# existing sheets) and records its summary row on "总计", ahead of the
# existing 2022-Q3 summary row.
#
# Approach (chosen so the resulting sheetId / r:id numbering matches a
# real-Excel "duplicate then edit" workflow):
#   1. Rename the existing "2022-Q3" sheet to "2022-Q4" - it keeps its
#      original sheetId/r:id slot.
#   2. Duplicate it right after itself; the duplicate becomes the "new"
#      "2022-Q3" sheet and keeps the original Q3 fund data untouched.
#   3. Overwrite the data on the (renamed) "2022-Q4" sheet with the actual
#      Q4 numbers.
#   4. Insert the 2022-Q4 row into "总计", pushing the existing 2022-Q3 row
#      down.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: rename + duplicate -------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Name = "2022-Q4"

$q3.Copy($null, $q3)
$q3copy = $wb.Worksheets.Item($q3.Index + 1)
$q3copy.Name = "2022-Q3"

$q4 = $wb.Worksheets.Item("2022-Q4")
$total = $wb.Worksheets.Item("总计")

# --- 3: overwrite the 2022-Q4 sheet with the real Q4 numbers -------------
# Header row + the "index" column (A) reuse the same cell style as the
# "总计" sheet's header/index cells (style index 2 in the original file),
# so copy formats from there instead of re-building bold+border by hand.
$total.Range("B1").Copy() | Out-Null
$q4.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$total.Range("A2").Copy() | Out-Null
$q4.Range("A2:A7").PasteSpecial(-4122) | Out-Null

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# data rows - columns B, D, E, F, G carry values that look numeric (fund
# codes with leading zeros, "1.40"-style figures); format them as text
# first so the assignment keeps the original text instead of being
# coerced into a Number (which would also strip leading zeros / trailing
# zero decimals).
$rows = @(
    @("014263","鑫元长三角混合A","1.40","89.50","5.74","0.0804",3),
    @("005262","鑫元欣享灵活配置混合A","1.55","83.07","3.94","0.0611",7),
    @("009395","鑫元安鑫回报混合A","3.61","20.13","1.10","0.0397",7),
    @("005263","鑫元欣享灵活配置混合C","1.00","83.07","3.94","0.0394",7),
    @("014264","鑫元长三角混合C","0.22","89.50","5.74","0.0126",3),
    @("016259","鑫元安鑫回报混合C","0.00","20.13","1.10",0,7)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = 2 + $r
    $data = $rows[$r]

    $q4.Cells.Item($row, 1).Value = $r

    $q4.Cells.Item($row, 2).NumberFormat = "@"
    $q4.Cells.Item($row, 2).Value = $data[0]

    $q4.Cells.Item($row, 3).Value = $data[1]

    $q4.Cells.Item($row, 4).NumberFormat = "@"
    $q4.Cells.Item($row, 4).Value = $data[2]

    $q4.Cells.Item($row, 5).NumberFormat = "@"
    $q4.Cells.Item($row, 5).Value = $data[3]

    $q4.Cells.Item($row, 6).NumberFormat = "@"
    $q4.Cells.Item($row, 6).Value = $data[4]

    if ($r -eq 5) {
        $q4.Cells.Item($row, 7).Value = 0
    } else {
        $q4.Cells.Item($row, 7).NumberFormat = "@"
        $q4.Cells.Item($row, 7).Value = $data[5]
    }

    $q4.Cells.Item($row, 8).Value = $data[6]
}

# --- 4: insert the new summary row on "总计" ------------------------------
$total.Range("A2:D2").Copy() | Out-Null
$total.Range("A3").PasteSpecial(-4122) | Out-Null

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 6
$total.Cells.Item(3, 4).Value = 0.22

$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 4).Value = 0.23
